# Add 2022-Q4 data:
#  - "总计" (Total) sheet gains a new row2 for 2022-Q4, pushing the existing
#    2022-Q3 / 2022-Q1 rows down by one.
#  - A brand-new "2022-Q4" tab is inserted right after "总计", carrying the
#    fund-holding table for the new quarter. The previously-existing
#    "2022-Q3" and "2022-Q1" tabs keep their original data and simply shift
#    one position to the right.

$wb = $excel.ActiveWorkbook

$total = $wb.Worksheets.Item(1)   # "总计"
$q3    = $wb.Worksheets.Item(2)   # "2022-Q3" (current data, about to become 2022-Q4's base)
$q1    = $wb.Worksheets.Item(3)   # "2022-Q1" (untouched)

# ---------------------------------------------------------------------------
# 1. "总计" sheet: insert a new row 2, fill it with the 2022-Q4 summary, and
#    fix up the running index in column A for the rows that shifted down.
# ---------------------------------------------------------------------------
$total.Rows.Item(2).Insert()
$total.Range("B2:D2").ClearFormats()

$total.Range("A3").Copy()
$total.Range("A2").PasteSpecial(-4122)   # xlPasteFormats

$total.Range("A2").Value = 0
$total.Range("B2").Value = "2022-Q4"
$total.Range("C2").Value = 3
$total.Range("D2").Value = 0.26

$total.Range("A3").Value = 1
$total.Range("A4").Value = 2

# ---------------------------------------------------------------------------
# 2. Duplicate the current "2022-Q3" sheet into a new tab placed right before
#    it; the duplicate becomes "2022-Q4" (new data) while the original sheet
#    keeps its old data / name and simply ends up one tab further right.
# ---------------------------------------------------------------------------
$q3.Copy($q3)
$q4 = $wb.Worksheets.Item(2)
$q4.Name = "2022-Q4"

# Helper-free, explicit cell-by-cell update of the new "2022-Q4" figures.
# Fund-code column (B) and all the percentage/size columns (D:G) are stored
# as plain text in this workbook, so force text formatting before writing
# them, then clear the formatting again so no stray number format/style is
# left behind on the cell.

$q4.Range("B2:B3").NumberFormat = "@"
$q4.Range("D2:G4").NumberFormat = "@"

$q4.Range("B2").Value = "001628"
$q4.Range("C2").Value = "招商体育文化休闲股票A"
$q4.Range("D2").Value = "2.33"
$q4.Range("E2").Value = "93.03"
$q4.Range("F2").Value = "5.49"
$q4.Range("G2").Value = "0.1279"
$q4.Range("H2").Value = 1

$q4.Range("B3").Value = "513360"
$q4.Range("C3").Value = "博时中证全球中国教育主题ETF（QDII）"
$q4.Range("D3").Value = "4.89"
$q4.Range("E3").Value = "99.23"
$q4.Range("F3").Value = "2.37"
$q4.Range("G3").Value = "0.1159"
$q4.Range("H3").Value = 10

$q4.Range("D4").Value = "0.29"
$q4.Range("E4").Value = "93.03"
$q4.Range("F4").Value = "5.49"
$q4.Range("G4").Value = "0.0159"
$q4.Range("H4").Value = 1

$q4.Range("B2:B3").ClearFormats()
$q4.Range("D2:G4").ClearFormats()
